$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 4.078147666666667
$ws.Range("H2").Value = 12.234443
$ws.Range("I2").Value = 0.8119037611005604
$ws.Range("J2").Value = 0.8119037611005603
$ws.Range("M2").Value = 0.6327629999999999
$ws.Range("N2").Value = 1.898289
$ws.Range("O2").Value = 0.1382544270550543
$ws.Range("P2").Value = 0.1382544270550544
$ws.Range("Q2").Value = 2.580500952003
$ws.Range("R2").Value = 23.224508568027
$ws.Range("S2").Value = 0.1122492893148017
$ws.Range("T2").Value = 0.1122492893148017
$ws.Range("G3").Value = 4.078147666666667
$ws.Range("H3").Value = 12.234443
$ws.Range("I3").Value = 0.8119037611005604
$ws.Range("J3").Value = 0.8119037611005603
$ws.Range("O3").Value = 0.4765301499162115
$ws.Range("P3").Value = 0.4765301499162115
$ws.Range("Q3").Value = 8.894373451254779
$ws.Range("R3").Value = 80.04936106129301
$ws.Range("S3").Value = 0.386896620994786
$ws.Range("T3").Value = 0.3868966209947859
$ws.Range("G4").Value = 4.078147666666667
$ws.Range("H4").Value = 12.234443
$ws.Range("I4").Value = 0.8119037611005604
$ws.Range("J4").Value = 0.8119037611005603
$ws.Range("M4").Value = 1.444396333333334
$ws.Range("N4").Value = 4.333189000000001
$ws.Range("O4").Value = 0.3155908096798033
$ws.Range("P4").Value = 0.3155908096798033
$ws.Range("Q4").Value = 5.890461536525224
$ws.Range("R4").Value = 53.01415382872701
$ws.Range("S4").Value = 0.2562293653478034
$ws.Range("T4").Value = 0.2562293653478034
$ws.Range("G5").Value = 4.078147666666667
$ws.Range("H5").Value = 12.234443
$ws.Range("I5").Value = 0.8119037611005604
$ws.Range("J5").Value = 0.8119037611005603
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3186579999999999
$ws.Range("N5").Value = 0.9559739999999999
$ws.Range("O5").Value = 0.06962461334893082
$ws.Range("P5").Value = 0.06962461334893082
$ws.Range("Q5").Value = 1.299534379164667
$ws.Range("R5").Value = 11.695809412482
$ws.Range("S5").Value = 0.05652848544316921
$ws.Range("T5").Value = 0.05652848544316921
$ws.Range("I6").Value = 0.07049404087934076
$ws.Range("J6").Value = 0.07049404087934076
$ws.Range("M6").Value = 0.6327629999999999
$ws.Range("N6").Value = 1.898289
$ws.Range("O6").Value = 0.1382544270550543
$ws.Range("P6").Value = 0.1382544270550544
$ws.Range("Q6").Value = 0.2240535742229999
$ws.Range("R6").Value = 2.016482168007
$ws.Range("S6").Value = 0.009746113232568836
$ws.Range("T6").Value = 0.009746113232568838
$ws.Range("I7").Value = 0.07049404087934076
$ws.Range("J7").Value = 0.07049404087934076
$ws.Range("O7").Value = 0.4765301499162115
$ws.Range("P7").Value = 0.4765301499162115
$ws.Range("S7").Value = 0.03359253586843179
$ws.Range("T7").Value = 0.0335925358684318
$ws.Range("I8").Value = 0.07049404087934076
$ws.Range("J8").Value = 0.07049404087934076
$ws.Range("M8").Value = 1.444396333333334
$ws.Range("N8").Value = 4.333189000000001
$ws.Range("O8").Value = 0.3155908096798033
$ws.Range("P8").Value = 0.3155908096798033
$ws.Range("Q8").Value = 0.5114429274118889
$ws.Range("R8").Value = 4.602986346707
$ws.Range("S8").Value = 0.02224727143871231
$ws.Range("T8").Value = 0.02224727143871231
$ws.Range("I9").Value = 0.07049404087934076
$ws.Range("J9").Value = 0.07049404087934076
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3186579999999999
$ws.Range("N9").Value = 0.9559739999999999
$ws.Range("O9").Value = 0.06962461334893082
$ws.Range("P9").Value = 0.06962461334893082
$ws.Range("Q9").Value = 0.1128328676846666
$ws.Range("R9").Value = 1.015495809162
$ws.Range("S9").Value = 0.004908120339627823
$ws.Range("T9").Value = 0.004908120339627823
$ws.Range("G10").Value = 0.5878206666666667
$ws.Range("H10").Value = 1.763462
$ws.Range("I10").Value = 0.1170271037560039
$ws.Range("J10").Value = 0.1170271037560039
$ws.Range("M10").Value = 0.6327629999999999
$ws.Range("N10").Value = 1.898289
$ws.Range("O10").Value = 0.1382544270550543
$ws.Range("P10").Value = 0.1382544270550544
$ws.Range("Q10").Value = 0.3719511685019999
$ws.Range("R10").Value = 3.347560516518
$ws.Range("S10").Value = 0.01617951517969872
$ws.Range("T10").Value = 0.01617951517969873
$ws.Range("G11").Value = 0.5878206666666667
$ws.Range("H11").Value = 1.763462
$ws.Range("I11").Value = 0.1170271037560039
$ws.Range("J11").Value = 0.1170271037560039
$ws.Range("O11").Value = 0.4765301499162115
$ws.Range("P11").Value = 0.4765301499162115
$ws.Range("Q11").Value = 1.282027272929111
$ws.Range("R11").Value = 11.538245456362
$ws.Range("S11").Value = 0.05576694329710859
$ws.Range("T11").Value = 0.0557669432971086
$ws.Range("G12").Value = 0.5878206666666667
$ws.Range("H12").Value = 1.763462
$ws.Range("I12").Value = 0.1170271037560039
$ws.Range("J12").Value = 0.1170271037560039
$ws.Range("M12").Value = 1.444396333333334
$ws.Range("N12").Value = 4.333189000000001
$ws.Range("O12").Value = 0.3155908096798033
$ws.Range("P12").Value = 0.3155908096798033
$ws.Range("Q12").Value = 0.849046015590889
$ws.Range("R12").Value = 7.641414140318002
$ws.Range("S12").Value = 0.03693267842883963
$ws.Range("T12").Value = 0.03693267842883963
$ws.Range("G13").Value = 0.5878206666666667
$ws.Range("H13").Value = 1.763462
$ws.Range("I13").Value = 0.1170271037560039
$ws.Range("J13").Value = 0.1170271037560039
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.3186579999999999
$ws.Range("N13").Value = 0.9559739999999999
$ws.Range("O13").Value = 0.06962461334893082
$ws.Range("P13").Value = 0.06962461334893082
$ws.Range("Q13").Value = 0.1873137579986666
$ws.Range("R13").Value = 1.685823821988
$ws.Range("S13").Value = 0.008147966850356984
$ws.Range("T13").Value = 0.008147966850356984
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.002888666666666667
$ws.Range("H14").Value = 0.008666
$ws.Range("I14").Value = 0.0005750942640950189
$ws.Range("J14").Value = 0.0005750942640950189
$ws.Range("M14").Value = 0.6327629999999999
$ws.Range("N14").Value = 1.898289
$ws.Range("O14").Value = 0.1382544270550543
$ws.Range("P14").Value = 0.1382544270550544
$ws.Range("Q14").Value = 0.001827841386
$ws.Range("R14").Value = 0.016450572474
$ws.Range("S14").Value = 0.00007950932798510495
$ws.Range("T14").Value = 0.00007950932798510496
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.002888666666666667
$ws.Range("H15").Value = 0.008666
$ws.Range("I15").Value = 0.0005750942640950189
$ws.Range("J15").Value = 0.0005750942640950189
$ws.Range("O15").Value = 0.4765301499162115
$ws.Range("P15").Value = 0.4765301499162115
$ws.Range("Q15").Value = 0.006300134818444444
$ws.Range("R15").Value = 0.056701213366
$ws.Range("S15").Value = 0.0002740497558851527
$ws.Range("T15").Value = 0.0002740497558851527
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.002888666666666667
$ws.Range("H16").Value = 0.008666
$ws.Range("I16").Value = 0.0005750942640950189
$ws.Range("J16").Value = 0.0005750942640950189
$ws.Range("M16").Value = 1.444396333333334
$ws.Range("N16").Value = 4.333189000000001
$ws.Range("O16").Value = 0.3155908096798033
$ws.Range("P16").Value = 0.3155908096798033
$ws.Range("Q16").Value = 0.004172379541555556
$ws.Range("R16").Value = 0.03755141587400001
$ws.Range("S16").Value = 0.0001814944644479576
$ws.Range("T16").Value = 0.0001814944644479576
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.002888666666666667
$ws.Range("H17").Value = 0.008666
$ws.Range("I17").Value = 0.0005750942640950189
$ws.Range("J17").Value = 0.0005750942640950189
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.3186579999999999
$ws.Range("N17").Value = 0.9559739999999999
$ws.Range("O17").Value = 0.06962461334893082
$ws.Range("P17").Value = 0.06962461334893082
$ws.Range("Q17").Value = 0.0009204967426666665
$ws.Range("R17").Value = 0.008284470683999999
$ws.Range("S17").Value = 0.0000400407157768036
$ws.Range("T17").Value = 0.0000400407157768036
